$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replicate the formatting of the last existing data rows (49-51 pattern: A=date style,
# C/D/E=currency-ish, F=integer) onto the two new rows by copying row 50/51 formats down.
$ws.Range("A50:F50").Copy()
$ws.Range("A52:F52").PasteSpecial(-4122)

$ws.Range("A51:F51").Copy()
$ws.Range("A53:F53").PasteSpecial(-4122)

# Row 52: 2025-11-26 (serial 45987), 四方坪站 (shared string index 4)
$ws.Cells.Item(52, 1).Value = 45987
$ws.Cells.Item(52, 2).Value = "四方坪站"
$ws.Cells.Item(52, 3).Value = 7978.38
$ws.Cells.Item(52, 4).Value = 7108.77
$ws.Cells.Item(52, 5).Value = 2656.3
$ws.Cells.Item(52, 6).Value = 342

# Row 53: 2025-11-26 (serial 45987), 高岭站 (shared string index 5)
$ws.Cells.Item(53, 1).Value = 45987
$ws.Cells.Item(53, 2).Value = "高岭站"
$ws.Cells.Item(53, 3).Value = 5264.22
$ws.Cells.Item(53, 4).Value = 4509.73
$ws.Cells.Item(53, 5).Value = 1349.22
$ws.Cells.Item(53, 6).Value = 189

# Update the selection to match the post-edit state (active cell G52)
$ws.Range("G52").Select()
